$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert a new row at position 8 (old row 8 -> 9, old row 9 -> 10) ---
$ws.Rows("8").Insert()

# --- Step 2: Copy formats from row 7 onto the newly inserted (blank) row 8 ---
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: Re-create the merges for row 8 that mirror row 7's layout ---
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# --- Step 4: Force the "text" number format (numFmtId 49) on the columns that
#     the source workbook re-formatted as Text before typing the new values in
#     (this mirrors the styles.xml numFmtId 0 -> 49 change for styles 8/9/12). ---
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("Q7").NumberFormat = "@"

$ws.Range("C8:G8").NumberFormat = "@"
$ws.Range("H8:K8").NumberFormat = "@"
$ws.Range("N8:O8").NumberFormat = "@"
$ws.Range("Q8").NumberFormat = "@"

# --- Step 5: Fill in row 7 (first item line) ---
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "CONTROLOC 40MG 14 GASTRORESISTANT TAB"
$ws.Range("H7").Value = "0:0"
$ws.Range("N7").Value = "188.00"
$ws.Range("Q7").Value = "1:0"

# L7 / P7 keep their original numeric-looking number formats but still hold
# literal text values, so flip to Text, assign, then restore the format.
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "188.0000"
$ws.Range("P7").NumberFormat = "0.00"

# --- Step 6: Fill in row 8 (second item line) ---
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "URSOFALK 250MG 20 CAPS."
$ws.Range("H8").Value = "0:0"
$ws.Range("N8").Value = "122.00"
$ws.Range("Q8").Value = "1:0"

$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "1"
$ws.Range("L8").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "122.0000"
$ws.Range("P8").NumberFormat = "0.00"

# --- Step 7: Row 9 (was row 8 pre-insert) now carries the running total ---
$ws.Range("P9").Value = 310

# --- Step 8: Row 10 (was row 9 pre-insert) footer timestamp update ---
$ws.Range("A10").Value = "Thursday, 17 July, 2025 9:13 AM"

Write-Output "edit applied"
